$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived metrics for rows 2-26 (columns G-T) per commit "update scripts wuth new tpm".
# Each entry: row number, then a map of column letter -> new numeric value.
$updates = @(
    @{ Row = 2; Cells = @{ "G" = 31.90834366666667; "H" = 95.725031; "I" = 0.1125536485145784; "J" = 0.1157863270269485; "M" = 0.6625783333333333; "N" = 1.987735; "O" = 0.2722635610602984; "P" = 0.280021914495151; "Q" = 21.14177716608722; "R" = 190.275994494785; "S" = 0.03064425715490829; "T" = 0.03242270896644777 } },
    @{ Row = 3; Cells = @{ "G" = 31.90834366666667; "H" = 95.725031; "I" = 0.1125536485145784; "J" = 0.1157863270269485; "O" = 0.1775534789114854; "P" = 0.1826129978482772; "Q" = 13.78736130384645; "R" = 124.086251734618; "S" = 0.01998429185794394; "T" = 0.02114408828823208 } },
    @{ Row = 4; Cells = @{ "G" = 31.90834366666667; "H" = 95.725031; "I" = 0.1125536485145784; "J" = 0.1157863270269485; "M" = 0.5835723333333332; "N" = 1.750717; "O" = 0.2397987884847841; "P" = 0.2466320339880352; "Q" = 18.62082656635855; "R" = 167.587439097227; "S" = 0.02699022855333813; "T" = 0.02855661734266014 } },
    @{ Row = 5; Cells = @{ "G" = 31.90834366666667; "H" = 95.725031; "I" = 0.1125536485145784; "J" = 0.1157863270269485; "M" = 0.202277; "N" = 0.404554; "O" = 0.08311870999996575; "P" = 0.0569914931299551; "Q" = 6.454324031862333; "R" = 38.72594419117399; "S" = 0.00935531407032132; "T" = 0.006598835661299073 } },
    @{ Row = 6; Cells = @{ "G" = 31.90834366666667; "H" = 95.725031; "I" = 0.1125536485145784; "J" = 0.1157863270269485; "K" = 3; "L" = 1; "M" = 0.5530713333333334; "N" = 1.659214; "O" = 0.2272654615434663; "P" = 0.2337415605385816; "Q" = 17.64759017618156; "R" = 158.828311585634; "S" = 0.02557955687806675; "T" = 0.0270640767683095 } },
    @{ Row = 7; Cells = @{ "I" = 0.2312918537506949; "J" = 0.2379348388122522; "M" = 0.6625783333333333; "N" = 1.987735; "O" = 0.2722635610602984; "P" = 0.280021914495151; "Q" = 43.44524497306777; "R" = 391.0072047576099; "S" = 0.06297234374640191; "T" = 0.06662696908930202 } },
    @{ Row = 8; Cells = @{ "I" = 0.2312918537506949; "J" = 0.2379348388122522; "O" = 0.1775534789114854; "P" = 0.1826129978482772; "S" = 0.04106667327732237; "T" = 0.043449994208052 } },
    @{ Row = 9; Cells = @{ "I" = 0.2312918537506949; "J" = 0.2379348388122522; "M" = 0.5835723333333332; "N" = 1.750717; "O" = 0.2397987884847841; "P" = 0.2466320339880352; "Q" = 38.26482350188244; "R" = 344.383411516942; "S" = 0.0554635063158165; "T" = 0.05868235325288108 } },
    @{ Row = 10; Cells = @{ "I" = 0.2312918537506949; "J" = 0.2379348388122522; "M" = 0.202277; "N" = 0.404554; "O" = 0.08311870999996575; "P" = 0.0569914931299551; "Q" = 13.26329790050067; "R" = 79.57978740300399; "S" = 0.0192246805172585; "T" = 0.01356026173154545 } },
    @{ Row = 11; Cells = @{ "I" = 0.2312918537506949; "J" = 0.2379348388122522; "K" = 3; "L" = 1; "M" = 0.5530713333333334; "N" = 1.659214; "O" = 0.2272654615434663; "P" = 0.2337415605385816; "Q" = 36.26487368424045; "R" = 326.383863158164; "S" = 0.05256464989389557; "T" = 0.0556152605304717 } },
    @{ Row = 12; Cells = @{ "G" = 85.57939900000001; "H" = 256.738197; "I" = 0.3018731932863474; "J" = 0.3105433607867011; "M" = 0.6625783333333333; "N" = 1.987735; "O" = 0.2722635610602984; "P" = 0.280021914495151; "Q" = 56.70305555708833; "R" = 510.327500013795; "S" = 0.08218907059278469; "T" = 0.08695894642125043 } },
    @{ Row = 13; Cells = @{ "G" = 85.57939900000001; "H" = 256.738197; "I" = 0.3018731932863474; "J" = 0.3105433607867011; "O" = 0.1775534789114854; "P" = 0.1826129978482772; "Q" = 36.97823072564067; "R" = 332.8040765307661; "S" = 0.05359863565811025; "T" = 0.05670925407513863 } },
    @{ Row = 14; Cells = @{ "G" = 85.57939900000001; "H" = 256.738197; "I" = 0.3018731932863474; "J" = 0.3105433607867011; "M" = 0.5835723333333332; "N" = 1.750717; "O" = 0.2397987884847841; "P" = 0.2466320339880352; "Q" = 49.94176955969433; "R" = 449.475926037249; "S" = 0.07238882602609917; "T" = 0.07658994071230434 } },
    @{ Row = 15; Cells = @{ "G" = 85.57939900000001; "H" = 256.738197; "I" = 0.3018731932863474; "J" = 0.3105433607867011; "M" = 0.202277; "N" = 0.404554; "O" = 0.08311870999996575; "P" = 0.0569914931299551; "Q" = 17.310744091523; "R" = 103.864464549138; "S" = 0.02509131040953151; "T" = 0.01769832981282844 } },
    @{ Row = 16; Cells = @{ "G" = 85.57939900000001; "H" = 256.738197; "I" = 0.3018731932863474; "J" = 0.3105433607867011; "K" = 3; "L" = 1; "M" = 0.5530713333333334; "N" = 1.659214; "O" = 0.2272654615434663; "P" = 0.2337415605385816; "Q" = 47.33151231079534; "R" = 425.983610797158; "S" = 0.06860535059982174; "T" = 0.07258688976517927 } },
    @{ Row = 17; Cells = @{ "G" = 23.7449455; "H" = 47.489891; "I" = 0.08375803763818537; "J" = 0.05744244731349463; "M" = 0.6625783333333333; "N" = 1.987735; "O" = 0.2722635610602984; "P" = 0.280021914495151; "Q" = 15.73288641448083; "R" = 94.39731848688498; "S" = 0.02280426159479485; "T" = 0.01608514407001161 } },
    @{ Row = 18; Cells = @{ "G" = 23.7449455; "H" = 47.489891; "I" = 0.08375803763818537; "J" = 0.05744244731349463; "O" = 0.1775534789114854; "P" = 0.1826129978482772; "Q" = 10.26001682094967; "R" = 61.560100925698; "S" = 0.01487153096945895; "T" = 0.01048973750765897 } },
    @{ Row = 19; Cells = @{ "G" = 23.7449455; "H" = 47.489891; "I" = 0.08375803763818537; "J" = 0.05744244731349463; "M" = 0.5835723333333332; "N" = 1.750717; "O" = 0.2397987884847841; "P" = 0.2466320339880352; "Q" = 13.85689325030783; "R" = 83.14135950184699; "S" = 0.0200850759514998; "T" = 0.01416714761817773 } },
    @{ Row = 20; Cells = @{ "G" = 23.7449455; "H" = 47.489891; "I" = 0.08375803763818537; "J" = 0.05744244731349463; "M" = 0.202277; "N" = 0.404554; "O" = 0.08311870999996575; "P" = 0.0569914931299551; "Q" = 4.8030563409035; "R" = 19.212225363614; "S" = 0.006961860040614546; "T" = 0.003273730841434837 } },
    @{ Row = 21; Cells = @{ "G" = 23.7449455; "H" = 47.489891; "I" = 0.08375803763818537; "J" = 0.05744244731349463; "K" = 3; "L" = 1; "M" = 0.5530713333333334; "N" = 1.659214; "O" = 0.2272654615434663; "P" = 0.2337415605385816; "Q" = 13.13264866761233; "R" = 78.79589200567401; "S" = 0.01903530908181722; "T" = 0.01342668727621149 } },
    @{ Row = 22; Cells = @{ "G" = 76.69186633333334; "H" = 230.075599; "I" = 0.270523266810194; "J" = 0.2782930260606035; "M" = 0.6625783333333333; "N" = 1.987735; "O" = 0.2722635610602984; "P" = 0.280021914495151; "Q" = 50.81436897536278; "R" = 457.329320778265; "S" = 0.07365362797140865; "T" = 0.07792814594813914 } },
    @{ Row = 23; Cells = @{ "G" = 76.69186633333334; "H" = 230.075599; "I" = 0.270523266810194; "J" = 0.2782930260606035; "O" = 0.1775534789114854; "P" = 0.1826129978482772; "Q" = 33.13799303561356; "R" = 298.241937320522; "S" = 0.04803234714864993; "T" = 0.05081992376919555 } },
    @{ Row = 24; Cells = @{ "G" = 76.69186633333334; "H" = 230.075599; "I" = 0.270523266810194; "J" = 0.2782930260606035; "M" = 0.5835723333333332; "N" = 1.750717; "O" = 0.2397987884847841; "P" = 0.2466320339880352; "Q" = 44.75525138383144; "R" = 402.797262454483; "S" = 0.06487115163803053; "T" = 0.06863597506201194 } },
    @{ Row = 25; Cells = @{ "G" = 76.69186633333334; "H" = 230.075599; "I" = 0.270523266810194; "J" = 0.2782930260606035; "M" = 0.202277; "N" = 0.404554; "O" = 0.08311870999996575; "P" = 0.0569914931299551; "Q" = 15.51300064630767; "R" = 93.078003877846; "S" = 0.02248554496223987; "T" = 0.0158603350828473 } },
    @{ Row = 26; Cells = @{ "G" = 76.69186633333334; "H" = 230.075599; "I" = 0.270523266810194; "J" = 0.2782930260606035; "K" = 3; "L" = 1; "M" = 0.5530713333333334; "N" = 1.659214; "O" = 0.2272654615434663; "P" = 0.2337415605385816; "Q" = 42.41607276879845; "R" = 381.744654919186; "S" = 0.06148059508986502; "T" = 0.07258688976517927 } }
)

foreach ($update in $updates) {
    $r = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $ws.Range("$col$r").Value = $update.Cells[$col]
    }
}

Write-Host "Applied all changes"